{"js": "// Append the missing student ID numbers after the two authors' names on\n// the title page:\n//   \"Chen Naveh\"        -> \"Chen Naveh - XXXXX\"\n//   \"Elad Wasserstein\"  -> \"Elad Wasserstein - 204499149\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  if (text === \"Chen Naveh\") {\n    para.insertText(\" - XXXXX\", Word.InsertLocation.end);\n  } else if (text === \"Elad Wasserstein\") {\n    para.insertText(\" - 204499149\", Word.InsertLocation.end);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Append the missing student ID numbers after the two authors' names on\n# the title page:\n#   \"Chen Naveh\"        -> \"Chen Naveh - XXXXX\"\n#   \"Elad Wasserstein\"  -> \"Elad Wasserstein - 204499149\"\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($text -eq \"Chen Naveh\") {\n        $insertRange = $p.Range\n        $insertRange.MoveEnd(1, -1) | Out-Null  # 1 = wdCharacter; shrink past the paragraph mark\n        $insertRange.InsertAfter(\" - XXXXX\")\n    }\n    elseif ($text -eq \"Elad Wasserstein\") {\n        $insertRange = $p.Range\n        $insertRange.MoveEnd(1, -1) | Out-Null\n        $insertRange.InsertAfter(\" - 204499149\")\n    }\n}\n"}
